# Gantt Chart update
# - Tracking column (F) gets "Done" / "Partially Done" status for rows 6-9
# - Date ranges in column C updated for rows 6-8 (UI/UX Designing, Frontend, Backend)
# - Column B widened; new column G added with its own width
# - Header row's G1 cell picks up the same header formatting as the rest of row 1
# - View zoomed out slightly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Dates" column for rows 6-8 ---
$ws.Range("C6").Value = "1-4June"
$ws.Range("C7").Value = "5-10June"
$ws.Range("C8").Value = "7-11June"

# --- Populate the "Tracking" column for rows 6-9 ---
$ws.Range("F6").Value = "Done"
$ws.Range("F7").Value = "Partially Done"
$ws.Range("F8").Value = "Partially Done"
$ws.Range("F9").Value = "Partially Done"

# New tracking cells pick up the same body font/style already used on D6 (Arial, no fill)
$ws.Range("D6").Copy() | Out-Null
$ws.Range("F6:F9").PasteSpecial(-4122) | Out-Null

# --- Header formatting: extend the header band into column G ---
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial(-4122) | Out-Null

# A blank but styled cell below the new Tracking entries (matches row formatting)
$ws.Range("D6").Copy() | Out-Null
$ws.Range("G7").PasteSpecial(-4122) | Out-Null

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 62.8
$ws.Columns.Item(7).ColumnWidth = 17.15

# --- View state ---
$excel.ActiveWindow.Zoom = 61
$ws.Range("G22").Select() | Out-Null
